# Update "RAF-generation" sheet values to latest 4.0 figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RAF-generation")

$ws.Range("B2").Value  = 0.55
$ws.Range("B3").Value  = 0.8
$ws.Range("B4").Value  = 0.8
$ws.Range("B12").Value = 0.8
$ws.Range("B13").Value = 0.8
$ws.Range("B14").Value = 0.55
$ws.Range("B16").Value = 0.8
$ws.Range("B17").Value = 0.8
$ws.Range("B18").Value = 0.8
